$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3104.9
$ws.Range("I40").Value = 3899.7778
$ws.Range("J40").Value = 2454.5454
$ws.Range("K40").Value = 3899.7778
$ws.Range("L40").Value = 2454.5454
$ws.Range("M40").Value = -3724.7778
$ws.Range("N40").Value = -2804.5454
$ws.Range("H74").Value = 4164.6
$ws.Range("I74").Value = 3958.25
$ws.Range("K74").Value = 3958.25
$ws.Range("M74").Value = -3022.25
$ws.Range("H77").Value = 4164.6
$ws.Range("I77").Value = 3958.25
$ws.Range("K77").Value = 19791.25
$ws.Range("M77").Value = -15111.25
$ws.Range("H86").Value = 2995
$ws.Range("I86").Value = 2975
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 2975
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1852
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 2995
$ws.Range("I89").Value = 2975
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 14875
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -9259
$ws.Range("N89").Value = -26232
$ws.Range("H98").Value = 1503.7646
$ws.Range("I98").Value = 1616.3846
$ws.Range("J98").Value = 1137.75
$ws.Range("K98").Value = 1616.3846
$ws.Range("L98").Value = 1137.75
$ws.Range("M98").Value = -118.3846000000001
$ws.Range("N98").Value = -4133.75
$ws.Range("H122").Value = 1503.7646
$ws.Range("I122").Value = 1616.3846
$ws.Range("J122").Value = 1137.75
$ws.Range("K122").Value = 4849.1538
$ws.Range("L122").Value = 3413.25
$ws.Range("M122").Value = -2399.1538
$ws.Range("N122").Value = -8313.25
$ws.Range("H138").Value = 3800.5908
$ws.Range("I138").Value = 3046.5806
$ws.Range("J138").Value = 5598.615
$ws.Range("K138").Value = 9139.7418
$ws.Range("L138").Value = 16795.845
$ws.Range("M138").Value = -3999.7418
$ws.Range("N138").Value = -27075.845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 577.3333
$ws.Range("I4").Value = 616
$ws.Range("K4").Value = 616
$ws.Range("M4").Value = -500
$ws.Range("H45").Value = 2707.2307
$ws.Range("J45").Value = 2499.3333
$ws.Range("L45").Value = 2499.3333
$ws.Range("N45").Value = -3253.3333
$ws.Range("H63").Value = 2724.1667
$ws.Range("I63").Value = 2724.1667
$ws.Range("K63").Value = 2724.1667
$ws.Range("M63").Value = -2038.1667
$ws.Range("H66").Value = 2724.1667
$ws.Range("I66").Value = 2724.1667
$ws.Range("K66").Value = 13620.8335
$ws.Range("M66").Value = -10188.8335
$ws.Range("H132").Value = 1650.3
$ws.Range("I132").Value = 1700.75
$ws.Range("J132").Value = 1448.5
$ws.Range("K132").Value = 5102.25
$ws.Range("L132").Value = 4345.5
$ws.Range("M132").Value = -2572.25
$ws.Range("N132").Value = -9405.5
$ws.Range("H133").Value = 150000
$ws.Range("J133").Value = 150000
$ws.Range("L133").Value = 150000
$ws.Range("N133").Value = -160120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H105").Value = 4953.353
$ws.Range("I105").Value = 1051.3334
$ws.Range("K105").Value = 1051.3334
$ws.Range("M105").Value = 695.6666
$ws.Range("H134").Value = 2452.85
$ws.Range("I134").Value = 2129.8125
$ws.Range("K134").Value = 6389.4375
$ws.Range("M134").Value = -3854.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1603
$ws.Range("I5").Value = 878.2
$ws.Range("J5").Value = 2327.8
$ws.Range("K5").Value = 2634.6
$ws.Range("L5").Value = 6983.400000000001
$ws.Range("M5").Value = -2522.6
$ws.Range("N5").Value = -7207.400000000001
$ws.Range("H21").Value = 2931
$ws.Range("I21").Value = 240.66667
$ws.Range("J21").Value = 11002
$ws.Range("K21").Value = 722.00001
$ws.Range("L21").Value = 33006
$ws.Range("M21").Value = -549.00001
$ws.Range("N21").Value = -33352
$ws.Range("H92").Value = 343
$ws.Range("I92").Value = 395
$ws.Range("J92").Value = 308.33334
$ws.Range("K92").Value = 1185
$ws.Range("L92").Value = 925.0000200000001
$ws.Range("M92").Value = 63
$ws.Range("N92").Value = -3421.00002
$ws.Range("H97").Value = 8694
$ws.Range("J97").Value = 9750
$ws.Range("L97").Value = 29250
$ws.Range("N97").Value = -30242
$ws.Range("H122").Value = 349.66666
$ws.Range("I122").Value = 277.8
$ws.Range("K122").Value = 2500.2
$ws.Range("M122").Value = -50.20000000000027
$ws.Range("H131").Value = 990
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("H135").Value = 1603
$ws.Range("I135").Value = 878.2
$ws.Range("J135").Value = 2327.8
$ws.Range("K135").Value = 7903.8
$ws.Range("L135").Value = 20950.2
$ws.Range("M135").Value = -5368.8
$ws.Range("N135").Value = -26020.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3289927.8
$ws.Range("I11").Value = 2859141.8
$ws.Range("J11").Value = 3720713.5
$ws.Range("K11").Value = 2859141.8
$ws.Range("L11").Value = 3720713.5
$ws.Range("M11").Value = -2859002.8
$ws.Range("N11").Value = -3720991.5
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H141").Value = 29995
$ws.Range("J141").Value = 29995
$ws.Range("L141").Value = 29995
$ws.Range("N141").Value = -40355

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2656.125
$ws.Range("I7").Value = 2535.5715
$ws.Range("K7").Value = 2535.5715
$ws.Range("M7").Value = -2423.5715
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 100
$ws.Range("K20").Value = 100
$ws.Range("M20").Value = 126
$ws.Range("H30").Value = 2366.3333
$ws.Range("I30").Value = 2366.3333
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 2366.3333
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -2258.3333
$ws.Range("N30").ClearContents()
$ws.Range("H126").Value = 2656.125
$ws.Range("I126").Value = 2535.5715
$ws.Range("K126").Value = 7606.7145
$ws.Range("M126").Value = -5136.7145
$ws.Range("H131").Value = 21904.762
$ws.Range("J131").Value = 21904.762
$ws.Range("L131").Value = 21904.762
$ws.Range("N131").Value = -31984.762

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 201739.4
$ws.Range("I3").Value = 500851.5
$ws.Range("J3").Value = 2331.3333
$ws.Range("K3").Value = 500851.5
$ws.Range("L3").Value = 2331.3333
$ws.Range("M3").Value = -500737.5
$ws.Range("N3").Value = -2559.3333
$ws.Range("H140").Value = 29999
$ws.Range("J140").Value = 29999
$ws.Range("L140").Value = 29999
$ws.Range("N140").Value = -40359
